# "add functions to db"
#
# - stickers sheet: add an "Answer" column (C) with reply texts, fix the
#   file.id values for "Привет"/"Пока" (they were swapped), widen column C.
# - add a new "users" sheet (after "stickers", becomes the active tab) with
#   an id/name/sex/grade header row and one seed row (id=123).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "stickers" (existing first sheet) - add Answer column + fix values
# ---------------------------------------------------------------------
$stickers = $wb.Worksheets.Item(1)

$stickers.Range("C1").Value = "Answer"
$stickers.Range("C2").Value = "Здравствуйте!"
$stickers.Range("C3").Value = "До свидания"

# the sticker file.id values were swapped between "Привет" and "Пока"
$stickers.Range("B2").Value = "CAACAgIAAxkBAAPwYgOYLTqiqDRaosBXwhdlM-YCNGgAAkQSAAKat6FL-O5ub_QpkuUjBA"
$stickers.Range("B3").Value = "CAACAgIAAxkBAAIBBWIDoFuXpZ1fjuFcFjJ3Rnv_QhkpAALREgACq03JSwiCbNDIBAFQIwQ"

# widen the new column (raw stored width 26 -> compensate for the
# character-width/pixel rounding the host applies)
$stickers.Columns.Item(3).ColumnWidth = 25.166666666666668

$stickers.Range("B11").Select()

# ---------------------------------------------------------------------
# 2. new "users" sheet, placed right after "stickers"
# ---------------------------------------------------------------------
$users = $wb.Worksheets.Add($null, $stickers)
$users.Name = "users"

$users.Columns.Item(1).ColumnWidth = 12.736979166666666
$users.Columns.Item(2).ColumnWidth = 18.451822916666668
$users.Columns.Item(3).ColumnWidth = 3.1666666666666665
$users.Columns.Item(4).ColumnWidth = 4.736979166666667
$users.Columns.Item(5).ColumnWidth = 100.73697916666667

$users.Range("A1").Value = "id"
$users.Range("B1").Value = "name"
$users.Range("C1").Value = "sex"
$users.Range("D1").Value = "grade"

$users.Range("A2").Value = 123

$users.Range("A2").Select()
